$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new price strings look numeric to Excel's auto-detection
# (e.g. "0.9983"), which would otherwise silently convert the cell to a
# Number and mangle values like "1.0000" -> 1. Force those to remain plain
# text (matching the workbook's original text-cell storage) by briefly
# switching to a text number format, then restore the original style so
# the cell's formatting is left untouched.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Price column (D) updates
$ws.Range('D2').Value = '29.398.17'
$ws.Range('D3').Value = '1.847.08'
$ws.Range('D12').Value = '1.846.53'
$ws.Range('D17').Value = '2.106.20'
$ws.Range('D19').Value = '29.435.81'
$ws.Range('D38').Value = '1.267.31'
$ws.Range('D44').Value = '2.011.28'
Set-TextValue 'D4' '0.9983'
Set-TextValue 'D5' '240.49'
Set-TextValue 'D6' '0.6266'
Set-TextValue 'D7' '1.0000'
Set-TextValue 'D9' '0.2906'
Set-TextValue 'D10' '24.45'
Set-TextValue 'D11' '0.07746'
Set-TextValue 'D13' '5.007'
Set-TextValue 'D14' '0.6805'
Set-TextValue 'D15' '0.00001053'
Set-TextValue 'D16' '82.21'
Set-TextValue 'D18' '6.179'
Set-TextValue 'D20' '229.79'
Set-TextValue 'D22' '0.9995'
Set-TextValue 'D23' '7.484'
Set-TextValue 'D24' '0.9998'
Set-TextValue 'D25' '159.18'
Set-TextValue 'D26' '0.1376'
Set-TextValue 'D27' '8.430'
Set-TextValue 'D28' '17.57'
Set-TextValue 'D29' '0.06505'
Set-TextValue 'D30' '1.414'
Set-TextValue 'D31' '1.479'
Set-TextValue 'D33' '4.100'
Set-TextValue 'D34' '1.834'
Set-TextValue 'D35' '1.143'
Set-TextValue 'D36' '0.6954'
Set-TextValue 'D37' '2.582'
Set-TextValue 'D39' '2.838'
Set-TextValue 'D40' '0.01835'
Set-TextValue 'D41' '6.766'
Set-TextValue 'D42' '0.9097'
Set-TextValue 'D43' '0.9996'
Set-TextValue 'D46' '66.39'
Set-TextValue 'D47' '1.747'
Set-TextValue 'D49' '7.081'
Set-TextValue 'D50' '0.1171'
Set-TextValue 'D51' '9.060'

# Volume(1h) column (E) updates
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -1.76%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('E17').Value = '  -3.69%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('E29').Value = '  +15.93%  '
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  +4.26%  '
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('E47').Value = '  +4.00%  '
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('E50').Value = '  +3.40%  '
$ws.Range('E51').Value = '  +0.45%  '
